# The workbook lists manga/comic convention events. The two data sheets
# "展览" (sheet 1) and "全部类型" (sheet 4) both held the same 4 event rows
# (rows 2-5). The edit drops the two oldest events (old rows 2 and 3,
# both dated 2024-10-01) and promotes the remaining two events (old rows
# 4 and 5) up to become the new rows 2 and 3, while also bumping their
# "想去人数" (interested-count) figures in column F.
#
# Sheets "演出" and "本地生活" only contain a header row and are left
# untouched.

$wb = $excel.ActiveWorkbook

function Update-EventSheet($ws) {
    # Remove the old row 2 (丽水·CCAC动漫游戏嘉年华) and, after the shift,
    # the old row 3 (丽水·熙梦动漫游戏展) which has now become row 2.
    # This leaves the former rows 4 and 5 as the new rows 2 and 3.
    $ws.Rows.Item(2).EntireRow.Delete()
    $ws.Rows.Item(2).EntireRow.Delete()

    # Renumber the sequence column (A) for the remaining two rows.
    $ws.Cells.Item(2, 1).Value = 1
    $ws.Cells.Item(3, 1).Value = 2

    # Update the "想去人数" (interested count) figures that changed.
    $ws.Cells.Item(2, 6).Value = 126
    $ws.Cells.Item(3, 6).Value = 46
}

Update-EventSheet $wb.Worksheets.Item(1)  # 展览
Update-EventSheet $wb.Worksheets.Item(4)  # 全部类型
